$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("N9").Value = 3929.35
$ws.Range("O9").Value = 2981.49

$ws.Range("N10").Value = 469660.57
$ws.Range("O10").Value = 462394.27

$ws.Range("M11").Value = 122297.85
$ws.Range("N11").Value = 71124.35
$ws.Range("O11").Value = 71124.35

$ws.Range("M12").Value = 34828.75
$ws.Range("N12").Value = 22681.66
$ws.Range("O12").Value = 22681.66

$ws.Range("N13").Value = 1899.72
$ws.Range("O13").Value = 1899.72

$ws.Range("N16").Value = 5671.85
$ws.Range("O16").Value = 5671.85

$ws.Range("N17").Value = 223.92
$ws.Range("O17").Value = 223.92

$ws.Range("N19").Value = 3612
$ws.Range("O19").Value = 3612

$ws.Range("K21").Value = 35692.96
$ws.Range("M21").Value = 5512.66

$ws.Range("K24").Value = 116970.61

$ws.Range("N26").Value = 71231.46
$ws.Range("O26").Value = 71231.46
